$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin "Price" figures are stored as display strings (often with thousands
# separators, e.g. "28.109.45"), not numeric values. Force text format first
# on the cells whose new price text would otherwise be re-interpreted as a
# number by Excel, so the stored cell type matches the source data.
$textCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D29",
    "D30",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = @{
    "D2" = "28.109.45"
    "E2" = "  -0.43%  "
    "D3" = "1.800.38"
    "E3" = "  -0.23%  "
    "D4" = "1.002"
    "E4" = "  -0.33%  "
    "D5" = "311.09"
    "E5" = "  -1.31%  "
    "D6" = "1.002"
    "E6" = "  -0.25%  "
    "D7" = "0.5095"
    "E7" = "  -2.69%  "
    "D8" = "0.3867"
    "E8" = "  +1.28%  "
    "D9" = "0.07736"
    "E9" = "  -2.66%  "
    "D10" = "1.098"
    "E10" = "  -0.04%  "
    "D11" = "40.83"
    "E11" = "  -2.01%  "
    "D12" = "6.331"
    "E12" = "  -0.28%  "
    "E13" = "  -0.34%  "
    "D14" = "20.25"
    "E14" = "  -2.03%  "
    "D15" = "1.804.11"
    "E15" = "  -0.08%  "
    "D16" = "7.270"
    "E16" = "  -1.12%  "
    "D17" = "92.00"
    "E17" = "  -0.61%  "
    "D18" = "0.00001073"
    "E18" = "  -1.65%  "
    "D19" = "0.06553"
    "E19" = "  -0.68%  "
    "E20" = "  -0.29%  "
    "D21" = "17.20"
    "E21" = "  -1.49%  "
    "D22" = "5.942"
    "E22" = "  -0.50%  "
    "D23" = "28.121.03"
    "E23" = "  -0.57%  "
    "D24" = "11.05"
    "E24" = "  -0.61%  "
    "D25" = "2.245"
    "E25" = "  +0.48%  "
    "D26" = "160.46"
    "E26" = "  +1.68%  "
    "D27" = "2.417"
    "E27" = "  -0.16%  "
    "D28" = "2.004.88"
    "E28" = "  -0.27%  "
    "D29" = "20.21"
    "E29" = "  -1.41%  "
    "D30" = "127.37"
    "E30" = "  +3.28%  "
    "E31" = "  -1.75%  "
    "D32" = "1.045"
    "E32" = "  -1.68%  "
    "D33" = "3.651"
    "E33" = "  -0.52%  "
    "D34" = "5.528"
    "E34" = "  -0.95%  "
    "D35" = "0.07017"
    "E35" = "  -3.04%  "
    "D36" = "9.043"
    "D37" = "0.02342"
    "E37" = "  +1.14%  "
    "D38" = "0.2163"
    "E38" = "  -0.79%  "
    "D39" = "5.021"
    "E39" = "  -0.52%  "
    "D40" = "11.46"
    "E40" = "  -5.94%  "
    "D41" = "0.6102"
    "E41" = "  -1.68%  "
    "E42" = "  -0.29%  "
    "D43" = "1.154"
    "E43" = "  -1.01%  "
    "D44" = "13.10"
    "E44" = "  -0.67%  "
    "B45" = "Decentraland"
    "C45" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D45" = "0.5903"
    "E45" = "  -2.56%  "
    "B46" = "WEMIXTOKEN"
    "C46" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D46" = "1.294"
    "E46" = "  -6.26%  "
    "D47" = "3.717"
    "E47" = "  -1.34%  "
    "D48" = "125.12"
    "E48" = "  -0.88%  "
    "D49" = "1.195"
    "E49" = "  -1.46%  "
    "D50" = "1.905"
    "E50" = "  -1.56%  "
    "D51" = "0.06734"
    "E51" = "  -1.40%  "
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
